$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 200 previously had no explicit style; the edited workbook applies the
# sheet's normal/default style (same as used by the rest of the data rows)
# to it. Re-write its contents (without touching formatting) so it picks
# up the worksheet's default column style.
$ws.Range("A200:C200").ClearContents()
$ws.Cells.Item(200, 1).Value = "2.ACE.ANE.0049 ANEL ELASTICO 04/03/2021 (FATI FERRAME)"
$ws.Cells.Item(200, 2).Value = 25
$ws.Cells.Item(200, 3).Value = "Unidade (un)"

# New rows appended at the bottom of the sheet.
$newRows = @(
    @("2.ACE.ANE.0050 ANEL 15/03/2021 (MAGELB)", 9.2, "Unidade (un)"),
    @("2.ACE.ANE.0051 ANEL ORING 25X32X4,8 mm 4-012-01-0408 20/04/2021 (MTM)", 17, "Unidade (un)"),
    @("2.ACE.ANE.0052 ANEL ORING 24,9X33,5X7 TYP IV 4-012-01-0608 20/04/2021 (MTM)", 20.7, "Unidade (un)"),
    @("2.ACE.ANE.0053 QORING DE VITON ALTA TEMPERATURA 29,74X3,53 mm 4-012-02-042 20/04/2021 (MTM)", 27.1, "Unidade (un)"),
    @("2.ACE.ANE.0054 QORING VITON VD ›29,74X3,53 FKM 70 COD HOMAG (70283) 4-012- 20/04/2021 (MTM)", 39.5, "Unidade (un)"),
    @("2.ACE.ANE.0055 ANEL VED NILOSRING MS0 23/06/2021 (NPX SOLUCOES)", 46.2, "Unidade (un)"),
    @("2.ACE.ANE.0056 ANEL TRAVA 35X1,5-FST-PHR DIN 472 23/06/2021 (NPX SOLUCOES)", 257.6, "Unidade (un)")
)

$row = 201
foreach ($data in $newRows) {
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $row = $row + 1
}

# Final row keeps the sheet's un-styled look, matching the old last row.
$ws.Cells.Item(208, 1).Value = "2.ACE.ANE.0057 ANEL ELASTICO E-12 13/12/2021 (CRV PARAFUSO)"
$ws.Cells.Item(208, 2).Value = 202.01
$ws.Cells.Item(208, 3).Value = "Unidade (un)"
